$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (it currently sits inside the
#    {s1f12} field placeholder earlier in the document).
$existingGoBack = $d.Bookmarks("_GoBack")
$existingGoBack.Delete()

# 2. Replace the run text "Copia Verdadera del Original" with the literal
#    placeholder "{o1}" (keeps the run's existing italic formatting).
$target = $d.Content
$target.Find.Execute("Copia Verdadera del Original", $true, $false, $false, $false, $false, $true, 1, $false, "{o1}", 2)

# 3. Locate the now-adjacent run that used to hold just a single space
#    (non-italic) and delete it entirely - it was merged/dropped in the
#    edit, its former position becomes where the new bookmark goes.
$afterNew = $d.Range($target.End, $target.End + 1)
if ($afterNew.Text -eq " ") {
    $afterNew.Delete()
}

# 4. Re-insert the "_GoBack" bookmark right after the "{o1}" run, where the
#    deleted space run used to live.
$bookmarkPoint = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)
